$d = $word.ActiveDocument
$lb = [char]11

# --- First remove the four paragraphs that are being dropped entirely ---
# (delete from the bottom up so indices of earlier paragraphs stay valid)
$d.Paragraphs.Item(13).Range.Delete()   # "המאמר גם טוען שאימון כזה..."
$d.Paragraphs.Item(12).Range.Delete()   # "אחרי שהגדרנו את ה-MDP..."
$d.Paragraphs.Item(11).Range.Delete()   # "המצב המתחלתי הוא רעש..."
$d.Paragraphs.Item(10).Range.Delete()   # "הפוליסי היא זו פונקצית..."

# --- Now update the remaining paragraphs' text (indices are stable again) ---

# Paragraph 1: title line + subtitle line (separated by a manual line break, no
# trailing break so a plain Range.Text assignment keeps the markup clean)
$d.Paragraphs.Item(1).Range.Text = "המאמר היומי של מייק - 25.12.24:" + $lb + "Vision language models are blind"

# Paragraph 2: intro (no trailing break)
$d.Paragraphs.Item(2).Range.Text = "מאמר נחמד הטוען שמודלי שפה ויזואליים הם די עיוורים כלומר אין להם סיכוי לעבור בדיקה אצל אופטומטריסט מורשה. הנה כמה עובדות על המבחנים הכושלים שלהם:"

# Paragraphs 3-7 each need new text that ends with a manual line break as the
# very last element of the run. Using Find/Replace (with the "^l" replacement
# code for the line break) keeps the serialized run clean, matching the
# target markup exactly (a plain Range.Text assignment would instead leave a
# spurious xml:space="preserve" attribute behind).

$r = $d.Paragraphs.Item(3).Range
$r.Find.Execute("קודם כל נשאלת השאלה למה צריך לאמן מודלי דיפוזיה גנרטיביים עם שיטות הלקוחות מעולם RL. הרי יש לנו שיטות סטנדרטיות יותר לאימון של מודלי דיפוזיה שהצליחו להביא לנו מודלים בעלי ביצועים מרשימים (בגנרוט תמונות מטקסט). אתם בטח יודעים שאימון מודלי דיפוזיה לגנרוט תמונות זה דבר לא זול ודורש לא מעט זמן ושימוש RL לאימון (או fine-tune) של מודלי דיפוזיה יכול לחסוך לנו זמן במקרים שאנו צריכים לאמן מודל דיפוזיה ייעודי (למשל לדומיין נישתי) ", $false, $false, $false, $false, $false, $true, 1, $false, "מודלי שפה ויזואליים או VLMs לא יכולים לקבוע באופן אמין האם שני קווים (או שני מעגלים) נחתכים, במיוחד כשהם קרובים זה לזה. הדיוק בזיהוי 0, 1 או 2 נקודות חיתוך בין שתי פונקציות לינאריות למקוטעין בעלות 2 מקטעים נע בין 47% ל-85%. באותה משימת שני המעגלים, המודלים מתפקדים טוב יותר (דיוק של 73-93%) אך עדיין רחוק מה-100% המצופה.^l", 2)

$r = $d.Paragraphs.Item(4).Range
$r.Find.Execute("אחת הדוגמאות למשימה כזו היא אימון מודל ליצירת תמונות מפרומפט (תיאור טקסטואלי) כאשר יש בידינו פונקציה המשערכת את התאמת התמונה לפרומפט. אתם כבר יכולים לנחש שפונקציה זו תשרת לנו בתור פונקצית תגמול (reward function). ", $false, $false, $false, $false, $false, $true, 1, $false, "מודלי שפה ויזואליים יכולים לזהות בצורה מושלמת מעגל ומילה בנפרד אך כאשר המעגל המילה נמצאת בתוך המעגל המודלים נוטים להתקשות בזיהוי איזו אות מוקפת במעגל.^l", 2)

$r = $d.Paragraphs.Item(5).Range
$r.Find.Execute("כבר הזכרתי שהמאמר משלב שיטה חדשה (יחסית) לאימון מודלי דיפוזיה הנקראת CM ושיטה זו (שהומצאה על ידי איליה סלוצקב ושות') מאפשרת גנרוט יותר מהיר של מודלי דיפוזיה גנרטיביים. בגדול מאוד שיטה זו מנסה לאמן מודל שאוכף עקביות בין התמונות המשוחזרות על ידי המודל מתמונות מורעשות עם עוצמות שונות רעש. כלומר לוקחים תמונה, מרעישים אותה עם רעש (בד״כ גאוסי) עם שונויות שונות ומאמנים מודל להחזיר את אותה התמונה הנקייה (עקביות לשמה). ", $false, $false, $false, $false, $false, $true, 1, $false, "מודלי ראייה-שפה יכולים לספור צורות במדויק, למשל, מעגלים , ריבועיים כאשר הם נפרדים ורחוקים זה מזה. עם זאת, כל המודלים מתקשים לספור מעגלים חותכים (כמו הלוגו האולימפי), ובאופן כללי, צורות בסיסיות שהן חופפות או מקוננות.^l", 2)

$r = $d.Paragraphs.Item(6).Range
$r.Find.Execute("למה השיטה הזו מאפשרת גנרוט יותר מהיר של תמונות? כי בגדול היא מאפשרת לגנרט תמונה נקייה מרעש באיטרציה אחת בלבד (ככה המודל מאומן). במציאות עושים את זה בכמה איטרציות (מספר קטן). מתחילים מרעש, מגנרטים את התמונה ממנו, מוסיפים פחות רעש לתמונה המגונרטת, מגנרטים מהתמונה המורעשת שוב וממשיכים ככה כמה איטרציות (עשרות בודדת). זה מאפשר לזרז את תהליך הגנרוט כי מודלי דיפוזיה סטנדרטיים צריכים מאות איטרציות בד״כ.", $false, $false, $false, $false, $false, $true, 1, $false, "בסידור ריבועים בצורה של רשת, אנו מגלים ש-VLMs נכשלים באופן מפתיע בספירת מספר השורות או העמודות ברשת, בין אם היא ריקה או מכילה טקסט. זה מפתיע בהתחשב בכך שהמודלים מתפקדים כל כך טוב (דיוק ≥ 90%) על הדאטהסט ב-DocVQA הכולל שאלות רבות עם טבלאות(אוברפיט כנראה).^l", 2)

$r = $d.Paragraphs.Item(7).Range
$r.Find.Execute("אוקיי, אחרי הקדמה ארוכה נעבור לתיאור של מה שעשו במאמר. המחברים הגדירו Markov Decision Process c או MDP המתאר תהליך גנרוט של תמונה (או כל דאטה אחר למעשה). כאמור פונקציה תגמול ניתנת לנו והיא מודדת מידת התאמה של התמונה המגונרטת לפרומפט. המאמר מגדיר:", $false, $false, $false, $false, $false, $true, 1, $false, "כאשר המודל מתבקש לעקוב אחר מסלולים צבעוניים במפת רכבת תחתית של עד 8 מסלולים וסך הכל 4 תחנות, VLMs לעתים קרובות נכשלים בזיהוי היכן מסלול מסתיים, כלומר, ומפגינים דיוק של 23% עד 50% .^l", 2)

# Paragraph 8: model comparison (no trailing break)
$d.Paragraphs.Item(8).Range.Text = "המודל GPT-4o עולה בביצועיו על Gemini-1.5 Pro ו-Claude-3 Sonnet ב-7 בנצ'מרקים מורכבים עבור VLMs אך מתפקד באופן משמעותי פחות טוב במשימות הנבחנות במאמר, שבהן Gemini-1.5 Pro ו-Sonnet-3.5 הם הטובים ביותר. כלומר, המאמר מגלה מגבלות מפתיעות של מודלי ראייה-שפה שלא נמדדו בבנצ'מרקים רגילים."

# Paragraph 9: closing joke
$d.Paragraphs.Item(9).Range.Text = "בקיצור אולי VLMs האלו צריכים משקפיים…"

# Paragraph 10 (was 14): updated arxiv link
$d.Paragraphs.Item(10).Range.Text = "https://arxiv.org/abs/2407.06581"
